$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1: "fourth" (same bold header style as A1:C1)
$ws.Range("D1").Value = "fourth"
$ws.Range("D1").Font.Bold = $true

# Give column D the same width as the other (custom-width) columns
$ws.Columns.Item(4).ColumnWidth = 8.29

# New cell D3: boolean TRUE
$ws.Range("D3").Value = $true
